# 动态规划 - coins in a line
# Restructure Sheet1: move the existing values/sums/dp table down, add a
# new "values" row, extend the data, add a title cell ("II"), a big notes
# cell with the DP derivation, and two small worked-example cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter       = -4108
$xlLeft         = -4131
$xlTop          = -4160

# ---------------------------------------------------------------------
# Row 8 (was row 2): column-index header, 0..4 now (one more column),
# reusing the existing italic/blue header style by copying its format
# down instead of re-deriving it (keeps the style table minimal/clean).
# ---------------------------------------------------------------------
$ws.Range("E2:H2").Copy() | Out-Null
$ws.Range("E8:I8").PasteSpecial($xlPasteFormats) | Out-Null

# Now that the old header's format has been copied forward, clear out
# the whole old D2:H5 block (content + format) so nothing old leaks
# into the rebuilt layout.
$ws.Range("D2:H5").Clear() | Out-Null

$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 4

# ---------------------------------------------------------------------
# Column C: big merged "II" title, centered both ways. Build the
# combined alignment on a scratch cell first (so the multi-cell range
# only ever gets ONE style write) then paste-format it across C4:C15.
# ---------------------------------------------------------------------
$scratch = $ws.Range("A1")
$scratch.HorizontalAlignment = $xlCenter
$scratch.VerticalAlignment = $xlCenter
$scratch.Copy() | Out-Null
$ws.Range("C4:C15").PasteSpecial($xlPasteFormats) | Out-Null
$scratch.Clear() | Out-Null

$ws.Range("C4").Value = "II"
$ws.Range("C4:C15").Merge() | Out-Null

# ---------------------------------------------------------------------
# K7:R20 big merged notes cell, left/top aligned with wrap text.
# Same scratch-cell technique to keep the style table clean.
# ---------------------------------------------------------------------
$scratch = $ws.Range("A1")
$scratch.HorizontalAlignment = $xlLeft
$scratch.VerticalAlignment = $xlTop
$scratch.WrapText = $true
$scratch.Copy() | Out-Null
$ws.Range("K7:R20").PasteSpecial($xlPasteFormats) | Out-Null
$scratch.Clear() | Out-Null

$bigText = @'
dp[i] = Math.max(sum[i] - dp[i + 1], sum[i] - dp[i + 2]);
dp[i] = Math.max(
        values[i] + Math.min(dp[i + 2], dp[i + 3]), 
        values[i] + values[i + 1] + Math.min(dp[i + 3], dp[i + 4])
     );
dp[i]=>从i到end的最大值
sum[i] - dp[i + 1]  => values[i] + min(dp[i + 2], dp[i + 3])
sum[i] - dp[i + 2]  => values[i] + values[i + 1] + min(dp[i + 3], dp[i + 4])
'@
$ws.Range("K7").Value = $bigText
$ws.Range("K7:R20").Merge() | Out-Null

# ---------------------------------------------------------------------
# Row 9: "values" row (was row 3), now bold, with an extra column.
# ---------------------------------------------------------------------
$ws.Range("D9").Value = "values"
$ws.Range("E9:H9").Font.Bold = $true
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 5

# ---------------------------------------------------------------------
# Row 10: "sums" row (was row 4), new values, extra column.
# ---------------------------------------------------------------------
$ws.Range("D10").Value = "sums"
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 7
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 0

# ---------------------------------------------------------------------
# Row 11: "dp" row (was row 5), new values, extra column.
# ---------------------------------------------------------------------
$ws.Range("D11").Value = "dp"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 0

# ---------------------------------------------------------------------
# E14:G14 and E15:G15: two small centered, merged worked-example cells.
# ---------------------------------------------------------------------
$ws.Range("E14:G14").HorizontalAlignment = $xlCenter
$ws.Range("E14").Value = "6=max(9-3,9-7)"
$ws.Range("E14:G14").Merge() | Out-Null

$ws.Range("E15:G15").HorizontalAlignment = $xlCenter
$ws.Range("E15").Value = "6=max(1+min(7,5),1+1+min(5,0))"
$ws.Range("E15:G15").Merge() | Out-Null

# ---------------------------------------------------------------------
# Column G width -> 11.5 (as stored in the raw column width units).
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 10.666666666667

# ---------------------------------------------------------------------
# Selection moves to J18.
# ---------------------------------------------------------------------
$ws.Range("J18").Select() | Out-Null
